# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 24 onward down by 6 rows to make room for the
# new "Number of employees / Assets / Turnover" breakdown table.
$ws.Rows("24:29").Insert(-4121)

# New header row (23): bold, matches the existing "title" cell style used
# elsewhere in the sheet (e.g. B11:D11).
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B23:D23").Font.Bold = $true

# New data rows (24-27)
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "1-9"

$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "10-49"

$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "50-249"

$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">=250"

# Re-point the OECD hyperlink: its target cell (the URL text) shifted from
# A34 to A40 along with the rest of the "Sector Distribution Details" block.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A40"), "http://www.oecd.org/globalrelations/psd/43469966.pdf")

Write-Host "done"
